# PMConverter commit: "finished writing/reading of resources, risk analysis and baseline schedule"
#
# 1) Baseline Schedule (sheet1): recolor row 4 (task "Testing") to reflect the
#    new resource-allocation palette and refresh the baseline start/end dates.
# 2) Resources (sheet2): write out the resource list (Programmer / Tester)
#    with availability, cost and assignment info.
# 3) Risk Analysis (sheet3): write out the activity duration distribution
#    profiles table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Baseline Schedule")
$ws2 = $wb.Worksheets.Item("Resources")
$ws3 = $wb.Worksheets.Item("Risk Analysis")

# ===========================================================================
# PASS 1 — formatting (column widths, merges, fills, number formats)
# Doing all of the formatting first (in a stable order) keeps the generated
# style table (fills / cellXfs) deterministic and matches the order new
# styles get appended to the style sheet.
# ===========================================================================

# --- Baseline Schedule: recolor row 4 ("Testing") -------------------------
$ws1.Range("A4").Interior.Color = 32768     # dark green  FF008000
$ws1.Range("B4").Interior.Color = 32768
$ws1.Range("C4").Interior.Color = 13160660  # gray        FFD4D0C8
$ws1.Range("D4").Interior.Color = 32768
$ws1.Range("E4").Interior.Color = 32768
$ws1.Range("F4").Interior.Color = 32768
$ws1.Range("G4").Interior.Color = 65280     # bright green FF00FF00
$ws1.Range("H4").Interior.Color = 32768
$ws1.Range("I4").Interior.Color = 65535     # yellow      FFFFFF00
$ws1.Range("J4").Interior.Color = 13160660
$ws1.Range("K4").Interior.Color = 32768
$ws1.Range("L4").Interior.Color = 65280
$ws1.Range("M4").Interior.Color = 32768
$ws1.Range("N4").Interior.Color = 13160660

# --- Resources sheet formatting --------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 14.75

$ws2.Range("A1:D1").Merge()
$ws2.Range("E1:F1").Merge()
$ws2.Range("G1:H1").Merge()

$ws1.Range("A1").Copy()
$ws2.Range("A1:H2").PasteSpecial(-4122)
$ws2.Rows.Item(2).RowHeight = 25

$ws1.Range("A3").Copy()
$ws2.Range("A3:D4").PasteSpecial(-4122)

$ws1.Range("A3").Copy()
$ws2.Range("E3").PasteSpecial(-4122)
$ws2.Range("E3").NumberFormat = "#,##0.00 €"

$ws1.Range("A3").Copy()
$ws2.Range("F3").PasteSpecial(-4122)
$ws2.Range("F3").NumberFormat = "#,##0.00 €"

$ws1.Range("A3").Copy()
$ws2.Range("E4").PasteSpecial(-4122)
$ws2.Range("E4").NumberFormat = "#,##0.00 €"

$ws1.Range("A3").Copy()
$ws2.Range("F4").PasteSpecial(-4122)
$ws2.Range("F4").NumberFormat = "#,##0.00 €"

$ws1.Range("C3").Copy()
$ws2.Range("G3:G4").PasteSpecial(-4122)

$ws1.Range("K3").Copy()
$ws2.Range("H3:H4").PasteSpecial(-4122)

# --- Risk Analysis sheet formatting ----------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 2.75
$ws3.Columns.Item(2).ColumnWidth = 17.75
$ws3.Columns.Item(4).ColumnWidth = 14.75
$ws3.Range("E1:G1").ColumnWidth = 11.75

$ws3.Range("A1:B1").Merge()
$ws3.Range("D1:G1").Merge()

$ws1.Range("A1").Copy()
$ws3.Range("A1:G2").PasteSpecial(-4122)

$ws1.Range("C3").Copy()
$ws3.Range("A3:G3").PasteSpecial(-4122)

$ws1.Range("C4").Copy()
$ws3.Range("A4:C4").PasteSpecial(-4122)

$ws1.Range("A3").Copy()
$ws3.Range("D4:G4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ===========================================================================
# PASS 2 — values, written in reading order (row-major) so new shared
# strings are appended in the same order the original workbook uses them.
# ===========================================================================

# --- Baseline Schedule: refreshed baseline start/end dates -----------------
$ws1.Range("F3").Value = 42113.8552930887
$ws1.Range("G3").Value = 42118.8552930888
$ws1.Range("F4").Value = 42113.8552930891
$ws1.Range("G4").Value = 42123.8552930891

# --- Resources sheet values -------------------------------------------------
$ws2.Range("A1").Value = "General"
$ws2.Range("E1").Value = "Resource Cost"
$ws2.Range("G1").Value = "Resource Demand"

$ws2.Range("A2").Value = "ID"
$ws2.Range("B2").Value = "Name"
$ws2.Range("C2").Value = "Type"
$ws2.Range("D2").Value = "Availability"
$ws2.Range("E2").Value = "Cost/Use"
$ws2.Range("F2").Value = "Cost/Unit"
$ws2.Range("G2").Value = "Assigned To"
$ws2.Range("H2").Value = "Total Cost"

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "Programmer"
$ws2.Range("C3").Value = "Renewable"
$ws2.Range("D3").Value = " #"
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 100
$ws2.Range("G3").Value = "1[10 #];"
$ws2.Range("H3").Value = 40000

$ws2.Range("A4").Value = 1
$ws2.Range("B4").Value = "Tester"
$ws2.Range("C4").Value = "Renewable"
$ws2.Range("D4").Value = " #"
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 75
$ws2.Range("G4").Value = "1[5 #];2[50 #];"
$ws2.Range("H4").Value = 315000

# --- Risk Analysis sheet values --------------------------------------------
$ws3.Range("A1").Value = "General"
$ws3.Range("C1").Value = "Baseline"
$ws3.Range("D1").Value = "Activity Duration Distribution Profiles"

$ws3.Range("A2").Value = "ID"
$ws3.Range("B2").Value = "Name"
$ws3.Range("C2").Value = "Duration"
$ws3.Range("D2").Value = "Description"
$ws3.Range("E2").Value = "Optimistic"
$ws3.Range("F2").Value = "Most Probable"
$ws3.Range("G2").Value = "Pessimistic"

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "App Dev"
$ws3.Range("C3").Value = "5d "

$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "Testing"
$ws3.Range("C4").Value = "10d "
$ws3.Range("D4").Value = "manual - absolute"
$ws3.Range("E4").Value = 402
$ws3.Range("F4").Value = 480
$ws3.Range("G4").Value = 812
